# "Minor fix to Datapath" -- evaluate_fit_DP.pptx
#
# 1. "11x4" -> "3x4"            (text-only edit on the left-hand block label)
# 2. "11x4" -> "3*3*4"          (text edit + widened box for the bottom-left label)
# 3. Nudge the DF / ADD1 / ADD2 big block-letter labels up/left a bit
# 4. Add three small "lv1" / "lv2" / "lv3" labels above the WR_interact_matrix area
# 5. Touch the presentation-level guides (best effort; empty guide list)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. "11x4" -> "3x4" (shape #29, id 95) ---------------------------------
$s.Shapes.Item(29).TextFrame.TextRange.Text = "3x4"

# --- 2. "11x4" -> "3*3*4" + widen box (shape #30, id 96) -------------------
$shp96 = $s.Shapes.Item(30)
$shp96.TextFrame.TextRange.Text = "3*3*4"
$shp96.Width = 57.11543467086616

# --- 3. Reposition DF / ADD1 / ADD2 labels ----------------------------------
$shpDF = $s.Shapes.Item(37)
$shpDF.Left = 169.68276220551144
$shpDF.Top = 6.271417422834646

$shpADD1 = $s.Shapes.Item(99)
$shpADD1.Left = 333.9767716535433
$shpADD1.Top = -9.015590551181102

$shpADD2 = $s.Shapes.Item(100)
$shpADD2.Left = 592.0191345582675
$shpADD2.Top = -10.280708861417322

# --- 4. Add lv1 / lv2 / lv3 labels ------------------------------------------
$lvWidth = 30.103307086614173
$lvHeight = 21.810944881889764

$lv1 = $s.Shapes.AddTextbox(1, 372.05889763779527, 43.321104062204746, $lvWidth, $lvHeight)
$lv1.Name = "文字方塊 212"
$lv1.TextFrame.WordWrap = -1
$lv1.TextFrame.AutoSize = 1
$lv1.Fill.Visible = 0
$lv1.TextFrame.TextRange.Text = "lv1"
$lv1.TextFrame.TextRange.Font.Name = "Times New Roman"
$lv1.TextFrame.TextRange.Font.NameComplexScript = "Times New Roman"
$lv1.TextFrame.TextRange.Font.Size = 12

$lv2 = $s.Shapes.AddTextbox(1, 424.4955905511811, 44.581968503937006, $lvWidth, $lvHeight)
$lv2.Name = "文字方塊 213"
$lv2.TextFrame.WordWrap = -1
$lv2.TextFrame.AutoSize = 1
$lv2.Fill.Visible = 0
$lv2.TextFrame.TextRange.Text = "lv2"
$lv2.TextFrame.TextRange.Font.Name = "Times New Roman"
$lv2.TextFrame.TextRange.Font.NameComplexScript = "Times New Roman"
$lv2.TextFrame.TextRange.Font.Size = 12

$lv3 = $s.Shapes.AddTextbox(1, 471.95338582677164, 44.38133858267717, $lvWidth, $lvHeight)
$lv3.Name = "文字方塊 214"
$lv3.TextFrame.WordWrap = -1
$lv3.TextFrame.AutoSize = 1
$lv3.Fill.Visible = 0
$lv3.TextFrame.TextRange.Text = "lv3"
$lv3.TextFrame.TextRange.Font.Name = "Times New Roman"
$lv3.TextFrame.TextRange.Font.NameComplexScript = "Times New Roman"
$lv3.TextFrame.TextRange.Font.Size = 12

# --- 5. Presentation guides (best effort - no visible content change) ------
try {
    $guides = $p.Guides
    $guides.Add(1, 3.0) | Out-Null
} catch {
    # Guides collection may be a no-op in this host; ignore.
}
